$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 3 ("Feedback" -> "Next steps")
# -----------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Title
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Next steps"

# Content placeholder body text
$body = $s3.Shapes.Item(3).TextFrame.TextRange

# First paragraph: replace "Your feedback is important to us!" with the
# certificate-of-attendance note (keeps the existing blank 2nd paragraph
# and the "Please complete the feedback form..." 3rd paragraph intact).
$para1 = $body.Paragraphs(1, 1)
$para1.Text = "A certificate of attendance will be issued automatically by LSHTM" + [char]0x2019 + "s short courses team (look out for an e-mail next week)."

# Append two new paragraphs (a blank line, then the new "Feel free..."
# line) right after the "Please complete the feedback form..." paragraph,
# before the trailing blank paragraph already at the end of the box.
$para3 = $body.Paragraphs(3, 1)
$null = $para3.InsertAfter("`r`rFeel free to contact us if you have any questions on the course material or about your modelling work!")

# -----------------------------------------------------------------------
# Slide 4 ("Further resources" -> "Recommended textbooks")
# -----------------------------------------------------------------------
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "Recommended textbooks"

# -----------------------------------------------------------------------
# Slide 6 ("Which models will we see in the course?" -> "Which models did
# we see in the course?")
# -----------------------------------------------------------------------
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "Which models did we see in the course?"
